$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 107
$ws1.Range("F4").Value = 405
$ws1.Range("F5").Value = 964
$ws1.Range("F6").Value = 5254
$ws1.Range("F7").Value = 441
$ws1.Range("F8").Value = 625
$ws1.Range("F9").Value = 909
$ws1.Range("F13").Value = 565
$ws1.Range("F17").Value = 1758
$ws1.Range("F18").Value = 1452
$ws1.Range("F19").Value = 824
$ws1.Range("F21").Value = 187
$ws1.Range("F23").Value = 508
$ws1.Range("F28").Value = 2548
$ws1.Range("F29").Value = 170
$ws1.Range("F30").Value = 97
$ws1.Range("F32").Value = 88
$ws1.Range("F34").Value = 251
$ws1.Range("F39").Value = 274
$ws1.Range("F40").Value = 636

$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F4").Value = 107
$ws2.Range("F5").Value = 964
$ws2.Range("F7").Value = 5254
$ws2.Range("F8").Value = 441
$ws2.Range("F9").Value = 625
$ws2.Range("F12").Value = 909
$ws2.Range("F18").Value = 565
$ws2.Range("F23").Value = 1758
$ws2.Range("F24").Value = 1452
$ws2.Range("F25").Value = 824
$ws2.Range("F26").Value = 187
$ws2.Range("F29").Value = 508
$ws2.Range("F33").Value = 2548
$ws2.Range("F34").Value = 170
$ws2.Range("F35").Value = 97
$ws2.Range("F36").Value = 88
$ws2.Range("F38").Value = 251
$ws2.Range("F42").Value = 274
